$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("BC1").Value = 0.84892197520185442
$ws.Range("A2").Value = 0.9214076188898499
$ws.Range("T2").Value = 0.85706582242758311
$ws.Range("H3").Value = 0.9215184330405326
$ws.Range("AD3").Value = 0.65361971362018489
$ws.Range("AN3").Value = 0.88995573451956056
$ws.Range("C5").Value = 0.86931044066131258
$ws.Range("Z5").Value = 0.99877772717152946
$ws.Range("AA5").Value = 0.89721737022191972
$ws.Range("AX5").Value = 0.91877870449489563
$ws.Range("D6").Value = 0.96828851279161354
$ws.Range("BP6").Value = 0.85618508514265246
$ws.Range("B7").Value = 0.97399756078742827
$ws.Range("T7").Value = 0.88943863170925663
$ws.Range("BP7").Value = 0.79201459615957903
$ws.Range("K8").Value = 0.73629270878278896
$ws.Range("I10").Value = 0.6734821408147913
$ws.Range("K10").Value = 0.69014462900448814
$ws.Range("L10").Value = 0.91003265659684751
$ws.Range("V10").Value = 0.76014910123359036
$ws.Range("BA10").Value = 0.93722139133882487
$ws.Range("M11").Value = 0.72215456977848436
$ws.Range("V12").Value = 0.93382542006045366
$ws.Range("AF12").Value = 0.65296125468578081
$ws.Range("L13").Value = 0.9113797286812908
$ws.Range("AU13").Value = 0.77732021129232565
$ws.Range("BM13").Value = 0.68243928989115821
$ws.Range("F14").Value = 0.65331949248378374
$ws.Range("AR14").Value = 0.82478370815194491
$ws.Range("D15").Value = 0.85745035192893249
$ws.Range("U15").Value = 0.99284325150558228
$ws.Range("U16").Value = 0.64728889352569086
$ws.Range("AB17").Value = 0.9907399758149753
$ws.Range("I18").Value = 0.94381600460354953
$ws.Range("AD18").Value = 0.8853777550643136
$ws.Range("A19").Value = 0.67834037941239855
$ws.Range("AU20").Value = 0.84011731970978731
$ws.Range("V21").Value = 0.80533074282435058
$ws.Range("BC21").Value = 0.94657977719914999
$ws.Range("I22").Value = 0.94171323209482627
$ws.Range("P22").Value = 0.87716038845066935
$ws.Range("T22").Value = 0.96400808557387596
$ws.Range("W22").Value = 0.83065961999859761
$ws.Range("AX22").Value = 0.89653553062758595
$ws.Range("X23").Value = 0.5357179668487233
$ws.Range("AH23").Value = 0.79706422435565372
$ws.Range("AS23").Value = 0.93516830303921017
$ws.Range("AR24").Value = 0.78173766274196066
$ws.Range("AY24").Value = 0.95509996620318227
$ws.Range("BP24").Value = 0.84094746361536243
$ws.Range("AH25").Value = 0.7694220463496324
$ws.Range("T26").Value = 0.90511342795354999
$ws.Range("BB27").Value = 0.83355062433352556
$ws.Range("AD28").Value = 0.9918568157891996
$ws.Range("AX28").Value = 0.97999905722573977
$ws.Range("BN28").Value = 0.74155609315959126
$ws.Range("N30").Value = 0.91429993199943538
$ws.Range("AF31").Value = 0.61075012192089551
$ws.Range("AO31").Value = 0.75223190350088542
$ws.Range("BH31").Value = 0.95192966312890748
$ws.Range("BG32").Value = 0.96955237486895696
$ws.Range("AH33").Value = 0.90599044243278049
$ws.Range("B34").Value = 0.8470181232686701
$ws.Range("AU35").Value = 0.8532695891729003
$ws.Range("BH35").Value = 0.63113149229564047
$ws.Range("AD36").Value = 0.89877826942234884
$ws.Range("AH36").Value = 0.93642599884652888
$ws.Range("BN36").Value = 0.98880335140018849
$ws.Range("H37").Value = 0.83443671216032833
$ws.Range("O37").Value = 0.60907155970167381
$ws.Range("AC37").Value = 0.82589720752407547
$ws.Range("AF37").Value = 0.98036363278298966
$ws.Range("N38").Value = 0.75752671727213627
$ws.Range("AF38").Value = 0.8857786333997083
$ws.Range("Y39").Value = 0.87936508436470806
$ws.Range("BA39").Value = 0.97918421218521745
$ws.Range("BB39").Value = 0.78578429450658405
$ws.Range("E40").Value = 0.95395960836827054
$ws.Range("AD40").Value = 0.91152595325432406
$ws.Range("AP40").Value = 0.68099779069214472
$ws.Range("BH42").Value = 0.80583919819852334
$ws.Range("M43").Value = 0.6910549470509737
$ws.Range("AO43").Value = 0.97726879186856042
$ws.Range("BC43").Value = 0.88302908923399881
$ws.Range("S44").Value = 0.94317411680165097
$ws.Range("S45").Value = 0.64736556033184378
$ws.Range("AG45").Value = 0.68172948737326489
$ws.Range("AQ45").Value = 0.89023666585321715
$ws.Range("AU45").Value = 0.83559094723206573
$ws.Range("AW45").Value = 0.8882023306583442
$ws.Range("F46").Value = 0.97511844814730009
$ws.Range("N46").Value = 0.65722853594477004
$ws.Range("X48").Value = 0.87299109618799053
$ws.Range("Y48").Value = 0.88185160381032546
$ws.Range("AP50").Value = 0.89268977847102249
$ws.Range("AW50").Value = 0.71971279505789021
$ws.Range("BC50").Value = 0.99314589079414417
$ws.Range("AV51").Value = 0.98508938309824323
$ws.Range("S52").Value = 0.97933868041102157
$ws.Range("X52").Value = 0.76412877547509472
$ws.Range("AJ52").Value = 0.91808217503350065
$ws.Range("AK52").Value = 0.88869387457149407
$ws.Range("BA52").Value = 0.91613478678839355
$ws.Range("A53").Value = 0.89931623710093778
$ws.Range("Z53").Value = 0.70935522835363041
$ws.Range("AI54").Value = 0.53551991171809332
$ws.Range("AC55").Value = 0.90155234351662961
$ws.Range("BG55").Value = 0.84611936907841634
$ws.Range("L56").Value = 0.70905830735678421
$ws.Range("W56").Value = 0.82997113978693848
$ws.Range("BB56").Value = 0.58445421931719088
$ws.Range("D57").Value = 0.93289074919402015
$ws.Range("AK57").Value = 0.75353701531730055
$ws.Range("BF57").Value = 0.72475629047610202
$ws.Range("O58").Value = 0.97399507759691495
$ws.Range("AE59").Value = 0.97045767969754337
$ws.Range("AW59").Value = 0.83142915537119055
$ws.Range("BJ59").Value = 0.7541904855406194
$ws.Range("N60").Value = 0.8785435366797083
$ws.Range("AS60").Value = 0.66739939751993038
$ws.Range("AT61").Value = 0.9333729181919227
$ws.Range("AY61").Value = 0.81766984836659518
$ws.Range("AL62").Value = 0.896214914489468
$ws.Range("BK62").Value = 0.89235702653929883
$ws.Range("F63").Value = 0.98600872639699255
$ws.Range("AT63").Value = 0.99123489020698208
$ws.Range("BF63").Value = 0.87838787545395736
$ws.Range("K64").Value = 0.72533955962475238
$ws.Range("Q64").Value = 0.8153965331121712
$ws.Range("AW64").Value = 0.9903047462337955
$ws.Range("BO65").Value = 0.90216146148867349
$ws.Range("AI66").Value = 0.81041887241748212
$ws.Range("BI66").Value = 0.96709123992446711
$ws.Range("Z67").Value = 0.63013367379260687
$ws.Range("AA67").Value = 0.97269174691881521
$ws.Range("AS68").Value = 0.57596955029275199
